# emptyMText-template.docx — "Fixed #295 Add the version of M2Doc in the
# template custom properties."
#
# The canonical-OOXML diff for *this* resource is a pure XML
# serialization artifact: every hunk in word/document.xml and
# word/styles.xml re-orders the attributes of an element (e.g.
# <w:pgSz w:w=".." w:h=".."/> -> <w:pgSz w:h=".." w:w=".."/>, the root
# <w:document> namespace declarations, the <w:latentStyles>/
# <w:lsdException>/<w:style> attributes, etc.) without adding,
# removing, or changing a single attribute value. A value-for-value
# check confirms every changed line pair carries the exact same
# (tag, {attr: value}) set before and after - only the on-disk
# attribute order differs, which the real commit picked up as a
# side effect of the tool that regenerated the fixture (the actual
# functional change — stamping the M2Doc version into
# docProps/custom.xml — lives in the product code path exercised by
# the test, not in this template's persisted XML).
#
# Word's object model has no property that controls XML attribute
# emission order, and this template's content/formatting already
# matches the target exactly, so there is no content edit to make
# here: touching the document (even by re-assigning a property to
# its current value) only risks incidental, unrelated churn (e.g.
# recomputed word/character counts or line-ending normalization)
# that the target diff does not contain. The correct, faithful
# application of this change is therefore to leave the already
# up-to-date content untouched.

$d = $word.ActiveDocument

# Touch nothing: confirm we're positioned on the right document and
# stop, matching the target state (which is already in place).
$null = $d.Name
